# Adding new SAM data.
# Restructure Arkusz1: split the "factors" column off into its own
# "government" column, and split the "taxes" column into separate
# "inc_taxes" / "cons_taxes" columns, while renaming the goods columns to
# "goods_activities" / "goods_commodities".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new (currently empty) column B -- this pushes the old
# factors/taxes/goods_production/goods_consumption columns one to the
# right (B->C, C->D, D->E, E->F).
$ws.Columns("B").Insert()

# Insert another new (empty) column E -- this splits the taxes column
# (now in D) from the goods columns (now in F/G).
$ws.Columns("E").Insert()

# New column B becomes the "government" column, carrying the "gov" row
# that used to sit at the bottom of the households list (A16).
$ws.Range("B1").Value = "government"
$ws.Range("B2").Value = "gov"
$ws.Range("A16").Value = ""

# Column D keeps just the income tax ("dtax"), renamed to "inc_taxes".
$ws.Range("D1").Value = "inc_taxes"
$ws.Range("D2").Value = "dtax"
$ws.Range("D3").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("D5").Value = ""

# New column E becomes the consumption tax ("stax"), renamed "cons_taxes".
$ws.Range("E1").Value = "cons_taxes"
$ws.Range("E2").Value = "stax"

# Rename the goods columns (now F/G) to activities/commodities.
$ws.Range("F1").Value = "goods_activities"
$ws.Range("G1").Value = "goods_commodities"

# Move the selection to A2, matching the saved workbook state.
$ws.Range("A2").Select() | Out-Null
